$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.416.41"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.641.14"
$ws.Range("E3").Value = "  +2.31%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "304.93"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.47%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.3732"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.96%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "52.34"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("E9").Value = "  -0.20%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "1.251"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.31%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.08117"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.06%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "22.77"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.45%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "6.600"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.00001269"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.61%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "7.287"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").Value = "1.631.80"
$ws.Range("E17").Value = "  +1.97%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "94.39"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.36%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.06888"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.83%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "18.12"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "6.510"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "23.415.47"
$ws.Range("E23").Value = "  +0.68%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "12.78"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.51%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "3.100"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.45%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "2.412"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.19%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "21.20"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "151.25"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.86%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "5.283"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.39%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "136.14"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.03%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "2.282"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.66%  "
$ws.Range("D32").Value = "1.812.62"
$ws.Range("E32").Value = "  +1.85%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "6.810"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.81%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.9504"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.79%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.02805"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("E36").Value = "  +1.06%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.2517"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.03%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.07227"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.63%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "6.117"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.21%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.08758"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("E41").Value = "  -1.70%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.7052"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "12.47"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "16.02"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.03%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.6513"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.86%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "2.328"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.9989"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.05%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "4.011"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +0.21%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "128.82"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.85%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "1.200"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
